$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(562, 6).Value = 27131
$ws.Cells.Item(564, 6).Value = 14324
$ws.Cells.Item(565, 6).Value = 29177
$ws.Cells.Item(565, 7).Value = 369
$ws.Cells.Item(566, 6).Value = 25774
$ws.Cells.Item(567, 6).Value = 23519
$ws.Cells.Item(568, 6).Value = 24115
$ws.Cells.Item(569, 6).Value = 32490
$ws.Cells.Item(570, 6).Value = 15233
$ws.Cells.Item(571, 6).Value = 15023
$ws.Cells.Item(571, 7).Value = 279
$ws.Cells.Item(572, 6).Value = 33563
$ws.Cells.Item(573, 6).Value = 27030
$ws.Cells.Item(574, 6).Value = 23459
$ws.Cells.Item(575, 6).Value = 26389
$ws.Cells.Item(576, 6).Value = 29016
$ws.Cells.Item(577, 6).Value = 14787
$ws.Cells.Item(577, 7).Value = 275
$ws.Cells.Item(578, 6).Value = 15099
$ws.Cells.Item(579, 6).Value = 32859
$ws.Cells.Item(580, 6).Value = 28895
$ws.Cells.Item(580, 7).Value = 518
$ws.Cells.Item(581, 6).Value = 27133
$ws.Cells.Item(582, 6).Value = 26144
$ws.Cells.Item(583, 6).Value = 29455
$ws.Cells.Item(584, 6).Value = 13268
$ws.Cells.Item(585, 6).Value = 14942
$ws.Cells.Item(586, 6).Value = 33878
$ws.Cells.Item(587, 6).Value = 28249
$ws.Cells.Item(588, 6).Value = 25400
$ws.Cells.Item(589, 6).Value = 26074
$ws.Cells.Item(590, 6).Value = 29043
$ws.Cells.Item(590, 7).Value = 579
$ws.Cells.Item(591, 6).Value = 14814
$ws.Cells.Item(592, 6).Value = 17982
$ws.Cells.Item(593, 6).Value = 37267
$ws.Cells.Item(593, 7).Value = 1195
$ws.Cells.Item(594, 6).Value = 29821
$ws.Cells.Item(594, 7).Value = 819
$ws.Cells.Item(595, 6).Value = 27281
$ws.Cells.Item(597, 6).Value = 29604
$ws.Cells.Item(597, 7).Value = 957
$ws.Cells.Item(598, 6).Value = 15497
$ws.Cells.Item(598, 7).Value = 710
$ws.Cells.Item(600, 6).Value = 40583
$ws.Cells.Item(600, 7).Value = 1676
$ws.Cells.Item(601, 6).Value = 31692
$ws.Cells.Item(601, 7).Value = 1334
$ws.Cells.Item(602, 6).Value = 30186
$ws.Cells.Item(602, 7).Value = 1302
$ws.Cells.Item(603, 6).Value = 32127
$ws.Cells.Item(603, 7).Value = 1534
$ws.Cells.Item(604, 6).Value = 30115
$ws.Cells.Item(604, 7).Value = 1538
$ws.Cells.Item(605, 6).Value = 14887
$ws.Cells.Item(605, 7).Value = 1036
$ws.Cells.Item(606, 6).Value = 14442
$ws.Cells.Item(606, 7).Value = 1276
$ws.Cells.Item(607, 6).Value = 10919
$ws.Cells.Item(607, 7).Value = 968
$ws.Cells.Item(608, 6).Value = 46174
$ws.Cells.Item(608, 7).Value = 2911
$ws.Cells.Item(609, 6).Value = 36436
$ws.Cells.Item(609, 7).Value = 2163
$ws.Cells.Item(610, 6).Value = 34008
$ws.Cells.Item(610, 7).Value = 1929
$ws.Cells.Item(611, 6).Value = 34124
$ws.Cells.Item(611, 7).Value = 2129
$ws.Cells.Item(612, 6).Value = 16302
$ws.Cells.Item(612, 7).Value = 1430
$ws.Cells.Item(613, 6).Value = 21629
$ws.Cells.Item(614, 6).Value = 47731
$ws.Cells.Item(614, 7).Value = 3334
$ws.Cells.Item(615, 6).Value = 36737
$ws.Cells.Item(615, 7).Value = 2359
$ws.Cells.Item(616, 6).Value = 38185
$ws.Cells.Item(616, 7).Value = 2569
$ws.Cells.Item(617, 6).Value = 38875
$ws.Cells.Item(617, 7).Value = 2599
$ws.Cells.Item(618, 6).Value = 37712
$ws.Cells.Item(618, 7).Value = 2657
$ws.Cells.Item(619, 6).Value = 17728
$ws.Cells.Item(619, 7).Value = 1889
$ws.Cells.Item(620, 6).Value = 25843
$ws.Cells.Item(620, 7).Value = 2394
$ws.Cells.Item(621, 6).Value = 55955
$ws.Cells.Item(621, 7).Value = 4108
$ws.Cells.Item(622, 6).Value = 41230
$ws.Cells.Item(622, 7).Value = 3000
$ws.Cells.Item(623, 6).Value = 14998
$ws.Cells.Item(623, 7).Value = 1567
$ws.Cells.Item(624, 6).Value = 51037
$ws.Cells.Item(624, 7).Value = 3928
$ws.Cells.Item(625, 6).Value = 43599
$ws.Cells.Item(625, 7).Value = 3547
$ws.Cells.Item(626, 6).Value = 20009
$ws.Cells.Item(626, 7).Value = 2082
$ws.Cells.Item(627, 6).Value = 33510
$ws.Cells.Item(627, 7).Value = 2711
$ws.Cells.Item(628, 6).Value = 64112
$ws.Cells.Item(628, 7).Value = 4177
$ws.Cells.Item(629, 6).Value = 45952
$ws.Cells.Item(629, 7).Value = 2910
$ws.Cells.Item(630, 6).Value = 46299
$ws.Cells.Item(630, 7).Value = 2936
$ws.Cells.Item(631, 6).Value = 41403
$ws.Cells.Item(631, 7).Value = 2716
$ws.Cells.Item(632, 6).Value = 43684
$ws.Cells.Item(632, 7).Value = 2612
$ws.Cells.Item(633, 6).Value = 23818
$ws.Cells.Item(633, 7).Value = 1918
$ws.Cells.Item(634, 6).Value = 45684
$ws.Cells.Item(634, 7).Value = 2116
$ws.Cells.Item(635, 6).Value = 79046
$ws.Cells.Item(635, 7).Value = 3462
$ws.Cells.Item(636, 6).Value = 45919
$ws.Cells.Item(636, 7).Value = 2083
$ws.Cells.Item(637, 6).Value = 30393
$ws.Cells.Item(637, 7).Value = 1388
